$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.728.92'
$ws.Range('E2').Value = '  -0.27%  '
$ws.Range('D3').Value = '1.634.52'
$ws.Range('E3').Value = '  -0.24%  '
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '215.32'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -0.06%  '
$ws.Range('E6').Value = '  -0.62%  '
$ws.Range('E7').Value = '  -0.09%  '
$ws.Range('E8').Value = '  -0.10%  '
$ws.Range('E9').Value = '  -1.14%  '
$ws.Range('E10').Value = '  -4.12%  '
$ws.Range('E11').Value = '  +1.11%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '4.23'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -0.46%  '
$ws.Range('B13').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C13').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D13').Value = '1.860.16'
$ws.Range('E13').Value = '  -0.27%  '
$ws.Range('B14').Value = 'WrappedEther'
$ws.Range('C14').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D14').Value = '1.634.33'
$ws.Range('E14').Value = '  +0.07%  '
$ws.Range('E15').Value = '  -0.99%  '
$ws.Range('D16').Value = '0.0₃0763'
$ws.Range('E16').Value = '  -0.16%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '62.71'
$ws.Range('D17').ClearFormats()
$ws.Range('D18').Value = '25.754.19'
$ws.Range('E18').Value = '  -0.31%  '
$ws.Range('E19').Value = '  -0.11%  '
$ws.Range('E20').Value = '  +1.48%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '193.65'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +0.77%  '
$ws.Range('E22').Value = '  +0.19%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.26'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +2.02%  '
$ws.Range('E24').Value = '  -0.09%  '
$ws.Range('E25').Value = '  +2.61%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '140.48'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +0.56%  '
$ws.Range('E27').Value = '  -1.25%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '6.88'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +1.01%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.50'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -0.23%  '
$ws.Range('E30').Value = '  -0.06%  '
$ws.Range('E31').Value = '  -0.37%  '
$ws.Range('E32').Value = '  +1.49%  '
$ws.Range('E33').Value = '  -0.12%  '
$ws.Range('E34').Value = '  +0.71%  '
$ws.Range('E35').Value = '  +0.03%  '
$ws.Range('E36').Value = '  -0.44%  '
$ws.Range('D37').Value = '1.122.34'
$ws.Range('E37').Value = '  -0.85%  '
$ws.Range('B38').Value = 'MXToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.52'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -1.75%  '
$ws.Range('B39').Value = 'ImmutableX'
$ws.Range('C39').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.547'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -1.81%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0156'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -0.94%  '
$ws.Range('E41').Value = '  +0.07%  '
$ws.Range('E42').Value = '  +2.14%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '99.62'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +0.75%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.801'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +0.16%  '
$ws.Range('D45').Value = '1.769.06'
$ws.Range('E45').Value = '  -0.45%  '
$ws.Range('D46').Value = '0.0₆0113'
$ws.Range('E46').Value = '  +0.68%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '55.01'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -0.94%  '
$ws.Range('E48').Value = '  -2.10%  '
$ws.Range('E49').Value = '  -0.46%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '7.56'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -3.15%  '
$ws.Range('B51').Value = 'SynthetixNetwork'
$ws.Range('C51').Value = 'https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.33'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +1.90%  '
